$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (z0bug.wt_1040): code "1040-20%" -> "1040-20%A"; add causale_pagamento_id (I2);
# wt_types (J2) "other" -> "ritenuta"
$ws.Range("B2").Value = "1040-20%A"
$ws.Range("I2").Value = "external.A"
$ws.Range("J2").Value = "ritenuta"

# Row 3 (z0bug.wt_1038): code "1040-23%" -> "1040-23%R"; add causale_pagamento_id (I3);
# wt_types (J3) "other" -> "ritenuta"
$ws.Range("B3").Value = "1040-23%R"
$ws.Range("I3").Value = "external.R"
$ws.Range("J3").Value = "ritenuta"

# Row 4 (z0bug.wt_enasarco_1): name "(A)" -> "(R)"; add causale_pagamento_id (I4)
$ws.Range("C4").Value = "Enasarco 17% su 50% (R)"
$ws.Range("I4").Value = "external.R"

# New row 5: z0bug.wt_1040-23A — duplicate formatting of G4 into G5 first (so the new
# row's journal_id cell carries the same style as the rest of the column), then fill in
# all the row's values.
$ws.Range("G4").Copy()
$ws.Range("G5").Insert(-4121)

$ws.Range("A5").Value = "z0bug.wt_1040-23A"
$ws.Range("B5").Value = "1040-23%A"
$ws.Range("C5").Value = "1040 – 23% su 100% (A)"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "z0bug.coa_153110"
$ws.Range("F5").Value = "z0bug.coa_260110"
$ws.Range("G5").Value = "z0bug.jou_misc"
$ws.Range("H5").Value = "account.account_payment_term_15days"
$ws.Range("I5").Value = "external.A"
$ws.Range("J5").Value = "ritenuta"
$ws.Range("K5").Value = 1

$ws.Range("A2").Select() | Out-Null
